$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(62, 8).Value = 7679
$ws.Cells.Item(62, 9).Value = 6859.5
$ws.Cells.Item(62, 11).Value = 6859.5
$ws.Cells.Item(62, 13).Value = -6235.5
$ws.Cells.Item(65, 8).Value = 7679
$ws.Cells.Item(65, 9).Value = 6859.5
$ws.Cells.Item(65, 11).Value = 34297.5
$ws.Cells.Item(65, 13).Value = -31177.5
$ws.Cells.Item(75, 8).Value = 45885.832
$ws.Cells.Item(75, 10).Value = 45885.832
$ws.Cells.Item(75, 12).Value = 45885.832
$ws.Cells.Item(75, 14).Value = -47757.832
$ws.Cells.Item(78, 8).Value = 45885.832
$ws.Cells.Item(78, 10).Value = 45885.832
$ws.Cells.Item(78, 12).Value = 137657.496
$ws.Cells.Item(78, 14).Value = -147017.496
$ws.Cells.Item(116, 8).Value = 7042.6665
$ws.Cells.Item(116, 9).Value = 2000
$ws.Cells.Item(116, 11).Value = 2000
$ws.Cells.Item(116, 13).Value = 1442
$ws.Cells.Item(138, 8).Value = 2294.1956
$ws.Cells.Item(138, 9).Value = 2910
$ws.Cells.Item(138, 11).Value = 8730
$ws.Cells.Item(138, 13).Value = -3590

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 894.5
$ws.Cells.Item(2, 9).Value = 894.5
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 894.5
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(5, 8).Value = 556.6
$ws.Cells.Item(5, 9).Value = 556.6
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 556.6
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(21, 8).Value = 800
$ws.Cells.Item(21, 9).Value = 800
$ws.Cells.Item(21, 11).Value = 800
$ws.Cells.Item(21, 13).Value = -426
$ws.Cells.Item(32, 8).Value = 10434.167
$ws.Cells.Item(32, 9).Value = 10434.167
$ws.Cells.Item(32, 11).Value = 10434.167
$ws.Cells.Item(32, 13).Value = -10147.167
$ws.Cells.Item(61, 8).Value = 3740.7856
$ws.Cells.Item(61, 9).Value = 1943.0834
$ws.Cells.Item(61, 11).Value = 1943.0834
$ws.Cells.Item(61, 13).Value = -1731.0834
$ws.Cells.Item(97, 8).Value = 774
$ws.Cells.Item(97, 9).Value = 700.3333
$ws.Cells.Item(97, 10).Value = 995
$ws.Cells.Item(97, 11).Value = 700.3333
$ws.Cells.Item(97, 12).Value = 995
$ws.Cells.Item(97, 13).Value = -204.3333
$ws.Cells.Item(97, 14).Value = -1987
$ws.Cells.Item(116, 8).Value = 894.5
$ws.Cells.Item(116, 9).Value = 894.5
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 894.5
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(133, 8).Value = 124568.37
$ws.Cells.Item(133, 10).Value = 124568.37
$ws.Cells.Item(133, 12).Value = 124568.37
$ws.Cells.Item(133, 14).Value = -129628.37
$ws.Cells.Item(135, 8).Value = 171141.72
$ws.Cells.Item(135, 10).Value = 171141.72
$ws.Cells.Item(135, 12).Value = 171141.72
$ws.Cells.Item(135, 14).Value = -181281.72
$ws.Cells.Item(136, 8).Value = 3740.7856
$ws.Cells.Item(136, 9).Value = 1943.0834
$ws.Cells.Item(136, 11).Value = 5829.2502
$ws.Cells.Item(136, 13).Value = -3279.2502

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 894.5
$ws.Cells.Item(3, 9).Value = 894.5
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 894.5
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(4, 8).Value = 556.6
$ws.Cells.Item(4, 9).Value = 556.6
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 556.6
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(5, 8).Value = 679
$ws.Cells.Item(5, 9).Value = 624
$ws.Cells.Item(5, 11).Value = 624
$ws.Cells.Item(5, 13).Value = -511
$ws.Cells.Item(94, 8).Value = 1458
$ws.Cells.Item(94, 9).Value = 1493.6154
$ws.Cells.Item(94, 11).Value = 1493.6154
$ws.Cells.Item(94, 13).Value = -1042.6154
$ws.Cells.Item(105, 8).Value = 4598.5
$ws.Cells.Item(105, 9).Value = 3918.2
$ws.Cells.Item(105, 11).Value = 3918.2
$ws.Cells.Item(105, 13).Value = -2171.2
$ws.Cells.Item(134, 8).Value = 1230.375
$ws.Cells.Item(134, 9).Value = 1230.375
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 3691.125
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(58, 8).Value = 1534.5238
$ws.Cells.Item(58, 9).Value = 1486.5
$ws.Cells.Item(58, 11).Value = 1486.5
$ws.Cells.Item(58, 13).Value = -1283.5
$ws.Cells.Item(105, 8).Value = 3545.8518
$ws.Cells.Item(105, 9).Value = 2836.2144
$ws.Cells.Item(105, 11).Value = 2836.2144
$ws.Cells.Item(105, 13).Value = -1089.2144
$ws.Cells.Item(136, 8).Value = 1534.5238
$ws.Cells.Item(136, 9).Value = 1486.5
$ws.Cells.Item(136, 11).Value = 4459.5
$ws.Cells.Item(136, 13).Value = -1909.5

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 6875440.5
$ws.Cells.Item(4, 9).Value = 6875440.5
$ws.Cells.Item(4, 11).Value = 20626321.5
$ws.Cells.Item(4, 13).Value = -20626209.5
$ws.Cells.Item(103, 8).Value = 1416.2858
$ws.Cells.Item(103, 9).Value = 880
$ws.Cells.Item(103, 11).Value = 2640
$ws.Cells.Item(103, 13).Value = -1761
$ws.Cells.Item(132, 8).Value = 4776.222
$ws.Cells.Item(132, 9).Value = 4997.8335
$ws.Cells.Item(132, 11).Value = 44980.5015
$ws.Cells.Item(132, 13).Value = -42450.5015
$ws.Cells.Item(134, 8).Value = 8291.6
$ws.Cells.Item(134, 9).Value = 1624.1666
$ws.Cells.Item(134, 10).Value = 18292.75
$ws.Cells.Item(134, 11).Value = 4872.4998
$ws.Cells.Item(134, 12).Value = 54878.25
$ws.Cells.Item(134, 13).Value = 197.5002000000004
$ws.Cells.Item(134, 14).Value = -65018.25
$ws.Cells.Item(140, 8).Value = 5671.1113
$ws.Cells.Item(140, 9).Value = 1006.9167
$ws.Cells.Item(140, 11).Value = 3020.7501
$ws.Cells.Item(140, 13).Value = 2159.2499

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(57, 8).Value = 45342
$ws.Cells.Item(57, 9).Value = 49685
$ws.Cells.Item(57, 10).Value = 40999
$ws.Cells.Item(57, 11).Value = 49685
$ws.Cells.Item(57, 12).Value = 40999
$ws.Cells.Item(57, 13).Value = -48865
$ws.Cells.Item(57, 14).Value = -42639
$ws.Cells.Item(102, 8).Value = 1186.5454
$ws.Cells.Item(102, 10).Value = 754.5
$ws.Cells.Item(102, 12).Value = 754.5
$ws.Cells.Item(102, 14).Value = -3998.5

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 7087.846
$ws.Cells.Item(7, 9).Value = 2899.8572
$ws.Cells.Item(7, 11).Value = 2899.8572
$ws.Cells.Item(7, 13).Value = -2787.8572
$ws.Cells.Item(82, 8).Value = 1504.75
$ws.Cells.Item(82, 9).Value = 1504.75
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 1504.75
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 1504.75
$ws.Cells.Item(85, 9).Value = 1504.75
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 1504.75
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 7087.846
$ws.Cells.Item(126, 9).Value = 2899.8572
$ws.Cells.Item(126, 11).Value = 8699.571599999999
$ws.Cells.Item(126, 13).Value = -6229.571599999999

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(45, 8).Value = 21202
$ws.Cells.Item(45, 10).Value = 21200
$ws.Cells.Item(45, 12).Value = 21200
$ws.Cells.Item(45, 14).Value = -22182
$ws.Cells.Item(132, 8).Value = 13500
$ws.Cells.Item(132, 9).Value = 12001
$ws.Cells.Item(132, 10).Value = 14999
$ws.Cells.Item(132, 11).Value = 36003
$ws.Cells.Item(132, 12).Value = 44997
$ws.Cells.Item(132, 13).Value = -33473
$ws.Cells.Item(132, 14).Value = -50057
$ws.Cells.Item(136, 8).Value = 2424.0667
$ws.Cells.Item(136, 9).Value = 2027.8462
$ws.Cells.Item(136, 11).Value = 6083.5386
$ws.Cells.Item(136, 13).Value = -3533.5386
